# Changed SHA1 generator flow
# - Adds new Settings entries for Acme System1 / SHA1 page URLs and log file locations
# - Inserts a new "Email" sheet (between Constants and Assets) with email
#   account / subject / body configuration used by the new flow
# - Adds new Assets entries used to route email notifications

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Settings sheet - new rows 7-15
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

$settings.Range("A7").Value = "AcmeSystem1LoginPageURL"
$settings.Range("B7").Value = "https://acme-test.uipath.com/login"

$settings.Range("A8").Value = "AcmeSystem1DashboardPageURL"
$settings.Range("B8").Value = "https://acme-test.uipath.com/"

$settings.Range("A9").Value = "AcmeSystem1WorkItemPageURL"
$settings.Range("B9").Value = "https://acme-test.uipath.com/work-items/"

$settings.Range("A10").Value = "SHA1LoginPageURL"
$settings.Range("B10").Value = "https://codebeautify.org/sha1-hash-generator/"

$settings.Range("A12").Value = "LogFileFolder"
$settings.Range("B12").Value = "D:\UiPath\Calculate Client Security Hash\Data\Logs\"

$settings.Range("A13").Value = "LogFileName"
$settings.Range("B13").Value = "Logs"

$settings.Range("A14").Value = "ExceptionScreenshotFolderPath"
$settings.Range("B14").Value = "D:\UiPath\Calculate Client Security Hash\Exceptions_Screenshots\"

$settings.Range("B15").Value = "Hi"

# ---------------------------------------------------------------------------
# 2. Insert a new "Email" worksheet between "Constants" and "Assets"
#
# NOTE: worksheet variables captured before an Add()/insert behave like
# *index*-based handles in this host, not stable object references - once a
# sheet is inserted before that index the variable silently starts pointing
# at the NEW sheet instead of the original one. So we look the "Assets"
# sheet up fresh (by name) both right before calling Add(), and again
# afterwards once the insertion has happened, rather than reusing a cached
# reference across the structural change.
# ---------------------------------------------------------------------------
$email = $wb.Worksheets.Add($wb.Worksheets.Item("Assets"))
$email.Name = "Email"

$assets = $wb.Worksheets.Item("Assets")

$email.Range("A1").Value = "Name"
$email.Range("B1").Value = "Value"
$email.Range("C1").Value = "Description"
$email.Range("A1:C1").Font.Bold = $true

$email.Range("A2").Value = "OutlookAccount"
$email.Range("B2").Value = "theafzalshaikh@gmail.com"

$email.Range("A3").Value = "Subject_Success"
$email.Range("B3").Value = "Success : Calculate Client Security Hash Process | "

$email.Range("A4").Value = "Subject_Failure"
$email.Range("B4").Value = "Failed : Calculate Client Security Hash Process | "

$email.Range("A5").Value = "Body_Success"
$email.Range("B5").Value = "Calculate client security hash process executed successfully." + [char]10 + "****This is a system generated email, do not reply to this email.****"

$email.Range("A6").Value = "Body_UnhandleException"
$email.Range("B6").Value = "Bot has encounter an unhandle exception." + [char]10 + "****This is a system generated email, do not reply to this email.****"

$email.Range("A7").Value = "Body_ACMESystemLoginPageNotResponding"
$b7 = "ACME System1 portal login page prompting username & password is not displayed." + [char]10 + "****It is a computer generated email, and you should not reply to it.****"
$email.Range("B7").Value = $b7
$b7italicStart = ("ACME System1 portal login page prompting username & password is not displayed." + [char]10).Length + 1
$email.Range("B7").Characters($b7italicStart, $b7.Length - $b7italicStart + 1).Font.Italic = $true

$email.Range("A8").Value = "Body_ACMESystemWrongCredentials"
$b8 = "ACME System1 wrong credentials." + [char]10 + "****It is a computer generated email, and you should not reply to it.****"
$email.Range("B8").Value = $b8
$b8italicStart = ("ACME System1 wrong credentials." + [char]10).Length + 1
$email.Range("B8").Characters($b8italicStart, $b8.Length - $b8italicStart + 1).Font.Italic = $true

$email.Range("A9").Value = "Body_ACMESystemHomePageNotResponding"
$b9 = "ACME System1 portal home page is not displayed." + [char]10 + "****It is a computer generated email, and you should not reply to it.****"
$email.Range("B9").Value = $b9
$b9italicStart = ("ACME System1 portal home page is not displayed." + [char]10).Length + 1
$email.Range("B9").Characters($b9italicStart, $b9.Length - $b9italicStart + 1).Font.Italic = $true

$email.Range("A10").Value = "Body_ACMESystemUnableToNavigateWorkItemsPage"
$b10 = "ACME System1 unable to navigate on work items page" + [char]10 + "****It is a computer generated email, and you should not reply to it.****"
$email.Range("B10").Value = $b10
$b10italicStart = ("ACME System1 unable to navigate on work items page" + [char]10).Length + 1
$email.Range("B10").Characters($b10italicStart, $b10.Length - $b10italicStart + 1).Font.Italic = $true

$email.Range("A11").Value = "Body_SHA1HomePageNotResponding"
$b11 = "SHA1 portal home page is not displayed." + [char]10 + "****It is a computer generated email, and you should not reply to it.****"
$email.Range("B11").Value = $b11
$b11italicStart = ("SHA1 portal home page is not displayed." + [char]10).Length + 1
$email.Range("B11").Characters($b11italicStart, $b11.Length - $b11italicStart + 1).Font.Italic = $true

# ---------------------------------------------------------------------------
# 3. Assets sheet - new rows 2-4
# ---------------------------------------------------------------------------
$assets.Range("A2").Value = "AcmeSystem1Credentials"
$assets.Range("B2").Value = "AcmeSystem1Credentials"

$assets.Range("A3").Value = "TOAddress"
$assets.Range("B3").Value = "CCSHTOAddress"

$assets.Range("A4").Value = "CCAddress"
$assets.Range("B4").Value = "CCSHCCAddress"

Write-Output "done"
